# correction des decimaux du prix unitaire d'un article
$wb = $excel.ActiveWorkbook

# --- Sheet "Facture_001" ---
$ws1 = $wb.Worksheets.Item("Facture_001")
$ws1.Range("D20").Value = 849999.5
$ws1.Range("G20").Value = 18
$ws1.Range("H20").Value = 849999.5
$ws1.Range("D21").Value = 24750.25
$ws1.Range("G21").Value = 18
$ws1.Range("H21").Value = 49500.5

# --- Sheet "Facture_002" ---
$ws2 = $wb.Worksheets.Item("Facture_002")
$ws2.Range("D20").Value = 649.75
$ws2.Range("G20").Value = 9
$ws2.Range("H20").Value = 32487.5

# --- Sheet "Facture_003_ERREURS" ---
$ws3 = $wb.Worksheets.Item("Facture_003_ERREURS")
$ws3.Range("G20").Value = 99

# --- Sheet "Avoir_001" ---
$ws4 = $wb.Worksheets.Item("Avoir_001")
$ws4.Range("D20").Value = -149999.99
$ws4.Range("G20").Value = 18
$ws4.Range("H20").Value = -149999.99

# --- Sheet "Facture_COMPLEXE" ---
$ws5 = $wb.Worksheets.Item("Facture_COMPLEXE")
$ws5.Range("D20").Value = 99999.75
$ws5.Range("G20").Value = 18
$ws5.Range("H20").Value = 99999.75

$ws5.Range("D21").Value = 499850.5
$ws5.Range("G21").Value = 18
$ws5.Range("H21").Value = 1499551.5

$ws5.Range("D22").Value = 74999.25
$ws5.Range("G22").Value = 9
$ws5.Range("H22").Value = 599994

$ws5.Range("D23").Value = 199950.75
$ws5.Range("G23").Value = 18
$ws5.Range("H23").Value = 199950.75

$ws5.Range("D24").Value = 349999
$ws5.Range("G24").ClearContents()
$ws5.Range("H24").Value = 699998

# --- New sheet "Test_TVA_Cases" appended at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws6 = $wb.Worksheets.Add($null, $lastSheet)
$ws6.Name = "Test_TVA_Cases"

$ws6.Range("A3").Value = "TEST001"
$ws6.Range("A5").Value = "CLI_TEST"
$ws6.Range("A6").Value = "1234567T"
$ws6.Range("A8").Value = "'45512"
$ws6.Range("A10").Value = "Test"
$ws6.Range("A11").Value = "CLIENT TEST TVA"
$ws6.Range("A18").Value = "cash"

$ws6.Range("B20").Value = "PROD_18"
$ws6.Range("C20").Value = "Produit avec TVA 18%"
$ws6.Range("D20").Value = 100.5
$ws6.Range("E20").Value = 1
$ws6.Range("F20").Value = "pcs"
$ws6.Range("G20").Value = 18
$ws6.Range("H20").Value = 100.5

$ws6.Range("B21").Value = "PROD_9"
$ws6.Range("C21").Value = "Produit avec TVA 9%"
$ws6.Range("D21").Value = 200.75
$ws6.Range("E21").Value = 1
$ws6.Range("F21").Value = "pcs"
$ws6.Range("G21").Value = 9
$ws6.Range("H21").Value = 200.75

$ws6.Range("B22").Value = "PROD_0"
$ws6.Range("C22").Value = "Produit avec TVA 0%"
$ws6.Range("D22").Value = 150
$ws6.Range("E22").Value = 1
$ws6.Range("F22").Value = "pcs"
$ws6.Range("G22").Value = 0
$ws6.Range("H22").Value = 150

$ws6.Range("B23").Value = "PROD_VIDE"
$ws6.Range("C23").Value = "Produit avec cellule TVA vide"
$ws6.Range("D23").Value = 300.25
$ws6.Range("E23").Value = 1
$ws6.Range("F23").Value = "pcs"
$ws6.Range("H23").Value = 300.25
